# Add working set of sequences
# Fill the distractor-image columns (G:N) with "N/A" for every trial row
# whose "cued" entry left those columns blank (i.e. only A:F were populated).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$rows = @(3, 5, 6, 10, 13, 16, 19, 20, 25, 28, 30, 32, 34, 35, 39, 41, 43, 48, 53, 59, 63, 66, 70, 72, 74, 77, 79, 85, 86, 87, 91, 92, 95, 97, 100, 102, 103, 104, 106, 107, 110, 111, 112, 119, 122, 128, 129, 135, 140, 146, 151, 152, 153, 155, 159, 164, 168, 173, 174, 176, 182, 186, 187, 193)

foreach ($r in $rows) {
    $ws.Range("G$r`:N$r").Value = "N/A"
}
